# Insert a new data row at row 88 (shifting existing rows 88-122 down to 89-123)
# and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(88).Insert()

$ws.Range("A88").Value = 8
$ws.Range("B88").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C88").Value = 'Coquimbo'
$ws.Range("D88").Value = 44704
$ws.Range("E88").Value = 4
$ws.Range("F88").Value = 'Fruta'
$ws.Range("G88").Value = 100109
$ws.Range("H88").Value = 'Uva'
$ws.Range("I88").Value = 100109001
$ws.Range("J88").Value = 'Uva'
$ws.Range("K88").Value = 'Red Globe'
$ws.Range("L88").Value = 'Primera'
$ws.Range("M88").Value = 300
$ws.Range("N88").Value = 9000
$ws.Range("O88").Value = 10000
$ws.Range("P88").Value = 9500
$ws.Range("Q88").Value = '$/bandeja 18 kilos'
$ws.Range("R88").Value = 'Provincia de Limarí'
$ws.Range("S88").Value = 528
$ws.Range("T88").Value = 18
